$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.452.32'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.223.20'
$ws.Range("E3").Value = '  -1.59%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.605'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '3.222.56'
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.57'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.410'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.59%  '
$ws.Range("D13").Value = '3.783.99'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").Value = '67.517.38'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").Value = '3.211.73'
$ws.Range("E18").Value = '  -2.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.75'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.05%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '396.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  -2.07%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.516'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("E26").Value = '  -3.33%  '
$ws.Range("E27").Value = '  -2.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.58'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("E35").Value = '  -2.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.20'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("E37").Value = '  -5.01%  '
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("E39").Value = '  -1.09%  '
$ws.Range("E40").Value = '  -4.36%  '
$ws.Range("E41").Value = '  -1.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.52'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.27%  '
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.96%  '
$ws.Range("D46").Value = '2.592.86'
$ws.Range("E46").Value = '  -2.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '332.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("E51").Value = '  -1.87%  '
